$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add two new columns (P1, Q1) ---
$ws.Cells.Item(1, 16).Value = 14   # P1
$ws.Cells.Item(1, 17).Value = 15   # Q1

# Copy the header formatting (bold, centered, bordered) from O1 onto the
# newly added P1:Q1 header cells so they match the rest of row 1.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-25: update existing values and add two new columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I{r} : 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K{r} : 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M{r} : 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O{r} : 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P{r} : new
    $ws.Cells.Item($r, 17).Value = 2   # Q{r} : new
}
